# Update TPM-derived NATMI LR-pair metrics (Hbegf-Egfr) on the active sheet
# to reflect the new TPM values used by the upstream scripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 7.723979
$ws.Range("H2").Value = 23.171937
$ws.Range("I2").Value = 0.471042132528101
$ws.Range("J2").Value = 0.471042132528101
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 3.168837049880334
$ws.Range("R2").Value = 28.519533448923
$ws.Range("S2").Value = 0.001648468408167481
$ws.Range("T2").Value = 0.001648468408167481

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 7.723979
$ws.Range("H3").Value = 23.171937
$ws.Range("I3").Value = 0.471042132528101
$ws.Range("J3").Value = 0.471042132528101
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("Q3").Value = 787.0593483903669
$ws.Range("R3").Value = 7083.534135513303
$ws.Range("S3").Value = 0.4094380527466375
$ws.Range("T3").Value = 0.4094380527466375

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 7.723979
$ws.Range("H4").Value = 23.171937
$ws.Range("I4").Value = 0.471042132528101
$ws.Range("J4").Value = 0.471042132528101
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("Q4").Value = 115.252170879713
$ws.Range("R4").Value = 1037.269537917417
$ws.Range("S4").Value = 0.05995561137329598
$ws.Range("T4").Value = 0.05995561137329598

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.2460132574367717
$ws.Range("J5").Value = 0.2460132574367717
$ws.Range("M5").Value = 0.4102596666666667
$ws.Range("N5").Value = 1.230779
$ws.Range("O5").Value = 0.003499619873322347
$ws.Range("P5").Value = 0.003499619873322347
$ws.Range("Q5").Value = 1.655002538187778
$ws.Range("R5").Value = 14.89502284369
$ws.Range("S5").Value = 0.0008609528848264931
$ws.Range("T5").Value = 0.000860952884826493

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.2460132574367717
$ws.Range("J6").Value = 0.2460132574367717
$ws.Range("O6").Value = 0.8692174743460166
$ws.Range("P6").Value = 0.8692174743460165
$ws.Range("S6").Value = 0.2138390222848271
$ws.Range("T6").Value = 0.2138390222848271

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.2460132574367717
$ws.Range("J7").Value = 0.2460132574367717
$ws.Range("N7").Value = 44.764041
$ws.Range("O7").Value = 0.1272829057806611
$ws.Range("P7").Value = 0.1272829057806611
$ws.Range("Q7").Value = 60.19326091405666
$ws.Range("R7").Value = 541.7393482265099
$ws.Range("S7").Value = 0.03131328226711815
$ws.Range("T7").Value = 0.03131328226711815

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.2829446100351274
$ws.Range("J8").Value = 0.2829446100351274
$ws.Range("M8").Value = 0.4102596666666667
$ws.Range("N8").Value = 1.230779
$ws.Range("O8").Value = 0.003499619873322347
$ws.Range("P8").Value = 0.003499619873322347
$ws.Range("Q8").Value = 1.903450458945444
$ws.Range("R8").Value = 17.131054130509
$ws.Range("S8").Value = 0.0009901985803283738
$ws.Range("T8").Value = 0.0009901985803283735

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.2829446100351274
$ws.Range("J9").Value = 0.2829446100351274
$ws.Range("O9").Value = 0.8692174743460166
$ws.Range("P9").Value = 0.8692174743460165
$ws.Range("S9").Value = 0.2459403993145521
$ws.Range("T9").Value = 0.245940399314552

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.2829446100351274
$ws.Range("J10").Value = 0.2829446100351274
$ws.Range("N10").Value = 44.764041
$ws.Range("O10").Value = 0.1272829057806611
$ws.Range("P10").Value = 0.1272829057806611
$ws.Range("Q10").Value = 69.22943467974565
$ws.Range("R10").Value = 623.0649121177109
$ws.Range("S10").Value = 0.03601401214024703
$ws.Range("T10").Value = 0.03601401214024703
